# Restore the "Rules" sheet's greeting text for the 06:00-11:00 rule row
# (E8) from "Good Morning" to "Good Morning1", per the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("E8").Value = "Good Morning1"
